# dados.xlsx — add "demanda_municipios" sheet + update selections on a few
# existing sheets, matching the commit "Add pyomo model and add parameter
# data in the dados.xlsx sheet".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) producao_maxima: selection moves from B1 to C10
# ---------------------------------------------------------------------
$wsProd = $wb.Worksheets.Item("producao_maxima")
$wsProd.Range("C10").Select()

# ---------------------------------------------------------------------
# 2) demanda_minima: selection becomes the range A2:A13, add a narrow
#    helper column E, and give the sheet an explicit portrait pageSetup.
# ---------------------------------------------------------------------
$wsDem = $wb.Worksheets.Item("demanda_minima")
$wsDem.Range("A2:A13").Select()
$wsDem.Columns.Item(5).ColumnWidth = 7.92
$wsDem.PageSetup.PaperSize = 9
$wsDem.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 3) distribuicao_maxima: selection moves from C2 to D3 (also loses
#    tabSelected once a later sheet is activated below).
# ---------------------------------------------------------------------
$wsDist = $wb.Worksheets.Item("distribuicao_maxima")
$wsDist.Range("D3").Select()

# ---------------------------------------------------------------------
# 4) Add the new "demanda_municipios" sheet at the end of the workbook
#    (after distribuicao_maxima) and populate it with the monthly
#    municipal-demand table.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsNew.Name = "demanda_municipios"

# Match the metric page-margin defaults used by every other sheet in this
# workbook (1.3cm / 2.0cm / 0.8cm, expressed in points for the COM API).
$wsNew.PageSetup.LeftMargin = 36.850393728
$wsNew.PageSetup.RightMargin = 36.850393728
$wsNew.PageSetup.TopMargin = 56.692913399999995
$wsNew.PageSetup.BottomMargin = 56.692913399999995
$wsNew.PageSetup.HeaderMargin = 22.67716464
$wsNew.PageSetup.FooterMargin = 22.67716464

$headerLabels = @(
    "Municipio","Recife","Jaboatão dos Guararapes","Cabo de Santo Agostinho",
    "Olinda","Paulista","Camaragibe","São Lourenço da Mata","Abreu e Lima",
    "Igarassu","Ipojuca"
)

$values = @(
    @(41275,41306,41334,41365,41395,41426,41456,41487,41518,41548,41579,41609),
    @(14400000,11800000,10800000,10800000,12100000,12800000,13100000,12900000,13000000,13800000,12300000,12200000),
    @(3600000,2500000,4700000,4200000,4100000,3700000,3800000,4300000,3800000,4800000,4500000,4400000),
    @(1800000,1500000,3000000,2600000,2500000,2900000,2800000,3000000,2700000,2900000,2400000,2500000),
    @(1500000,1400000,1500000,1400000,1400000,1400000,1700000,1600000,1300000,1400000,1500000,1900000),
    @(1300000,1400000,1200000,1100000,1400000,1300000,1500000,1700000,1300000,1500000,1100000,1100000),
    @(1200000,1200000,1300000,1100000,1200000,1200000,800000,1100000,1000000,1000000,1100000,1100000),
    @(700000,500000,350000,450000,500000,500000,750000,750000,750000,750000,750000,750000),
    @(300000,200000,300000,300000,200000,300000,300000,250000,500000,350000,350000,450000),
    @(200000,180000,200000,200000,210000,210000,240000,270000,260000,270000,280000,240000),
    @(130000,110000,170000,160000,170000,160000,170000,170000,160000,170000,160000,110000)
)

$data = New-Object 'object[,]' 11,13
for ($r = 0; $r -lt 11; $r++) {
    $data[$r,0] = $headerLabels[$r]
    for ($c = 0; $c -lt 12; $c++) {
        $data[$r,$c+1] = $values[$r][$c]
    }
}
$wsNew.Range("A1:M11").Value = $data

# Row 1 (B1:M1) holds month-end dates -> date number format.
$wsNew.Range("B1:M1").NumberFormat = "mmm-yy"
# Column A holds municipality names -> text number format (matches the
# rest of the workbook's "N." style text columns).
$wsNew.Range("A1:A11").NumberFormat = "@"

$wsNew.Columns.Item(1).ColumnWidth = 16.92

$wsNew.Range("B1").Select()
